# Update scripts with new TPM values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (Sending cluster): "MuSCs" -> "ECs" for rows 2-4
$ws.Range("A2").Value = "ECs"
$ws.Range("A3").Value = "ECs"
$ws.Range("A4").Value = "ECs"

# Row 2 (MuSCs-Camp-Fpr2-FAPs) updated values
$ws.Range("G2").Value = 0.05807666666666667
$ws.Range("H2").Value = 0.17423
$ws.Range("O2").Value = 0.2013489143128838
$ws.Range("P2").Value = 0.2013489143128838
$ws.Range("Q2").Value = 0.07304459173666666
$ws.Range("R2").Value = 0.65740132563
$ws.Range("S2").Value = 0.2013489143128838
$ws.Range("T2").Value = 0.2013489143128838

# Row 3 (MuSCs-Camp-Fpr2-MuSCs) updated values
$ws.Range("G3").Value = 0.05807666666666667
$ws.Range("H3").Value = 0.17423
$ws.Range("M3").Value = 0.246708
$ws.Range("N3").Value = 0.740124
$ws.Range("O3").Value = 0.03949536580856015
$ws.Range("P3").Value = 0.03949536580856015
$ws.Range("Q3").Value = 0.01432797828
$ws.Range("R3").Value = 0.12895180452
$ws.Range("S3").Value = 0.03949536580856015
$ws.Range("T3").Value = 0.03949536580856015

# Row 4 (MuSCs-Camp-Fpr2-Resolving-Mac) updated values
$ws.Range("G4").Value = 0.05807666666666667
$ws.Range("H4").Value = 0.17423
$ws.Range("M4").Value = 4.74207
$ws.Range("N4").Value = 14.22621
$ws.Range("O4").Value = 0.7591557198785561
$ws.Range("P4").Value = 0.759155719878556
$ws.Range("Q4").Value = 0.2754036187
$ws.Range("R4").Value = 2.4786325683
$ws.Range("S4").Value = 0.7591557198785561
$ws.Range("T4").Value = 0.759155719878556
